# Auto-generated script applying market-data refresh values
# to the Leve profit columns (H-N) across all profession sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 188.6
$ws.Range("J4").Value = 296.33334
$ws.Range("L4").Value = 296.33334
$ws.Range("N4").Value = -524.33334
$ws.Range("H12").Value = 409.6
$ws.Range("J12").Value = 217
$ws.Range("L12").Value = 217
$ws.Range("N12").Value = -557
$ws.Range("H106").Value = 3392.7693
$ws.Range("I106").Value = 3210.8
$ws.Range("J106").Value = 3999.3333
$ws.Range("K106").Value = 3210.8
$ws.Range("L106").Value = 3999.3333
$ws.Range("M106").Value = -2579.8
$ws.Range("N106").Value = -5261.3333
$ws.Range("H111").Value = 7328.5864
$ws.Range("I111").Value = 9758.223
$ws.Range("J111").Value = 3352.818
$ws.Range("K111").Value = 29274.669
$ws.Range("L111").Value = 10058.454
$ws.Range("M111").Value = -26207.669
$ws.Range("N111").Value = -16192.454
$ws.Range("H132").Value = 11508.491
$ws.Range("I132").Value = 2361.1277
$ws.Range("K132").Value = 7083.3831
$ws.Range("M132").Value = -4553.3831
$ws.Range("H133").Value = 39314.844
$ws.Range("J133").Value = 39314.844
$ws.Range("L133").Value = 39314.844
$ws.Range("N133").Value = -49434.844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1494
$ws.Range("J4").Value = 1989
$ws.Range("L4").Value = 1989
$ws.Range("N4").Value = -2221
$ws.Range("H6").Value = 2508248.5
$ws.Range("J6").Value = 9427
$ws.Range("L6").Value = 9427
$ws.Range("N6").Value = -9773
$ws.Range("H32").Value = 5608.963
$ws.Range("I32").Value = 5698.9863
$ws.Range("K32").Value = 5698.9863
$ws.Range("M32").Value = -5411.9863
$ws.Range("H74").Value = 3774.353
$ws.Range("I74").Value = 1001
$ws.Range("K74").Value = 1001
$ws.Range("M74").Value = -127
$ws.Range("H77").Value = 3774.353
$ws.Range("I77").Value = 1001
$ws.Range("K77").Value = 5005
$ws.Range("M77").Value = -637
$ws.Range("H97").Value = 1157.8572
$ws.Range("I97").Value = 1165.8
$ws.Range("K97").Value = 1165.8
$ws.Range("M97").Value = -669.8
$ws.Range("H122").Value = 2896.0444
$ws.Range("I122").Value = 2569.361
$ws.Range("J122").Value = 4202.778
$ws.Range("K122").Value = 7708.083
$ws.Range("L122").Value = 12608.334
$ws.Range("M122").Value = -5258.083
$ws.Range("N122").Value = -17508.334
$ws.Range("H132").Value = 2268.5903
$ws.Range("I132").Value = 2217.6353
$ws.Range("J132").Value = 2687.5557
$ws.Range("K132").Value = 6652.9059
$ws.Range("L132").Value = 8062.6671
$ws.Range("M132").Value = -4122.9059
$ws.Range("N132").Value = -13122.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1356.0714
$ws.Range("I94").Value = 976.1111
$ws.Range("J94").Value = 2040
$ws.Range("K94").Value = 976.1111
$ws.Range("L94").Value = 2040
$ws.Range("M94").Value = -525.1111
$ws.Range("N94").Value = -2942
$ws.Range("H99").Value = 2870.8948
$ws.Range("I99").Value = 2184.6155
$ws.Range("J99").Value = 4357.8335
$ws.Range("K99").Value = 2184.6155
$ws.Range("L99").Value = 4357.8335
$ws.Range("M99").Value = -686.6154999999999
$ws.Range("N99").Value = -7353.8335
$ws.Range("H134").Value = 6869.647
$ws.Range("I134").Value = 3255.6667
$ws.Range("K134").Value = 9767.000100000001
$ws.Range("M134").Value = -7232.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 278.9091
$ws.Range("I7").Value = 96.8125
$ws.Range("J7").Value = 764.5
$ws.Range("K7").Value = 96.8125
$ws.Range("L7").Value = 764.5
$ws.Range("M7").Value = 16.1875
$ws.Range("N7").Value = -990.5
$ws.Range("H31").Value = 1919.8918
$ws.Range("J31").Value = 5066.5
$ws.Range("L31").Value = 5066.5
$ws.Range("N31").Value = -5656.5
$ws.Range("H34").Value = 1919.8918
$ws.Range("J34").Value = 5066.5
$ws.Range("L34").Value = 5066.5
$ws.Range("N34").Value = -5470.5
$ws.Range("H99").Value = 5500
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4002
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 5500
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -14030
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2108081.2
$ws.Range("I132").Value = 2859990.5
$ws.Range("K132").Value = 8579971.5
$ws.Range("M132").Value = -8577441.5
$ws.Range("H134").Value = 4670.6113
$ws.Range("I134").Value = 1952.3334
$ws.Range("K134").Value = 5857.0002
$ws.Range("M134").Value = -3322.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2420
$ws.Range("I25").Value = 366.66666
$ws.Range("K25").Value = 1099.99998
$ws.Range("M25").Value = -930.9999800000001
$ws.Range("H30").Value = 2420
$ws.Range("I30").Value = 366.66666
$ws.Range("K30").Value = 1099.99998
$ws.Range("M30").Value = -997.9999800000001
$ws.Range("H100").Value = 10028
$ws.Range("J100").Value = 10028
$ws.Range("L100").Value = 30084
$ws.Range("N100").Value = -31706
$ws.Range("H121").Value = 4388.25
$ws.Range("I121").Value = 450
$ws.Range("J121").Value = 5701
$ws.Range("K121").Value = 1350
$ws.Range("L121").Value = 17103
$ws.Range("M121").Value = -40
$ws.Range("N121").Value = -19723
$ws.Range("H131").Value = 910456.2
$ws.Range("I131").Value = 1429616.9
$ws.Range("K131").Value = 4288850.699999999
$ws.Range("M131").Value = -4283810.699999999
$ws.Range("H136").Value = 3863.875
$ws.Range("I136").Value = 2574.7273
$ws.Range("K136").Value = 7724.1819
$ws.Range("M136").Value = -2624.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2399.5
$ws.Range("I113").Value = 2399.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2399.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -229.5
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 4904552
$ws.Range("I132").Value = 6063072
$ws.Range("J132").Value = 3120.4614
$ws.Range("K132").Value = 18189216
$ws.Range("L132").Value = 9361.3842
$ws.Range("M132").Value = -18186686
$ws.Range("N132").Value = -14421.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3707
$ws.Range("I7").Value = 2930.842
$ws.Range("K7").Value = 2930.842
$ws.Range("M7").Value = -2818.842
$ws.Range("H30").Value = 5150
$ws.Range("I30").Value = 5150
$ws.Range("K30").Value = 5150
$ws.Range("M30").Value = -5042
$ws.Range("H40").Value = 5373.7188
$ws.Range("I40").Value = 5318.346
$ws.Range("K40").Value = 5318.346
$ws.Range("M40").Value = -5182.346
$ws.Range("H56").Value = 13482.8
$ws.Range("I56").Value = 14978.5
$ws.Range("K56").Value = 14978.5
$ws.Range("M56").Value = -14287.5
$ws.Range("H122").Value = 3693.3333
$ws.Range("J122").Value = 5979.4
$ws.Range("L122").Value = 17938.2
$ws.Range("N122").Value = -22838.2
$ws.Range("H126").Value = 3707
$ws.Range("I126").Value = 2930.842
$ws.Range("K126").Value = 8792.526
$ws.Range("M126").Value = -6322.526
$ws.Range("H132").Value = 2600.4531
$ws.Range("I132").Value = 2488.6538
$ws.Range("K132").Value = 7465.9614
$ws.Range("M132").Value = -4935.9614
$ws.Range("H136").Value = 3998.5518
$ws.Range("I136").Value = 3236.6843
$ws.Range("K136").Value = 9710.052899999999
$ws.Range("M136").Value = -7160.052899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 833.05554
$ws.Range("I107").Value = 460.375
$ws.Range("J107").Value = 1131.2
$ws.Range("K107").Value = 1381.125
$ws.Range("L107").Value = 3393.6
$ws.Range("M107").Value = 538.875
$ws.Range("N107").Value = -7233.6
$ws.Range("H132").Value = 2077.2068
$ws.Range("I132").Value = 2003.5454
$ws.Range("J132").Value = 2308.7144
$ws.Range("K132").Value = 6010.6362
$ws.Range("L132").Value = 6926.1432
$ws.Range("M132").Value = -3480.6362
$ws.Range("N132").Value = -11986.1432
$ws.Range("H141").Value = 94643.7
$ws.Range("J141").Value = 93087.44500000001
$ws.Range("L141").Value = 93087.44500000001
$ws.Range("N141").Value = -103447.445
